$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct the existing row 235 close value ---
$ws.Range("F235").Value = 457.82

# --- Append new rows 236-238 with the same look (date format/border/bold) as column A's existing data cells ---
$ws.Range("A235").Copy()
$ws.Range("A236:A238").PasteSpecial(-4122)   # xlPasteFormats

$rows = @(
    @{ Row = 236; A = 45170.33333333334; B = "FX_IDC:USDKZT"; C = 457.82; D = 485.93; E = 454.52; F = 477.37; G = 0 },
    @{ Row = 237; A = 45201.375;         B = "FX_IDC:USDKZT"; C = 477.37; D = 480.77; E = 467.91; F = 468.27; G = 0 },
    @{ Row = 238; A = 45231.375;         B = "FX_IDC:USDKZT"; C = 468.27; D = 470.91; E = 461.68; F = 464.87; G = 0 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.A
    $ws.Range("B$row").Value = $r.B
    $ws.Range("C$row").Value = $r.C
    $ws.Range("D$row").Value = $r.D
    $ws.Range("E$row").Value = $r.E
    $ws.Range("F$row").Value = $r.F
    $ws.Range("G$row").Value = $r.G
}

Write-Output "Applied Kazakhstan_FX updates through row 238"
